$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($Range, $Value)
    $Range.NumberFormat = "@"
    $Range.Value = $Value
    $Range.Style = "Normal"
}

Set-TextCell ($ws.Range('D2')) '43.024.84'
Set-TextCell ($ws.Range('E2')) '  -5.27%  '
Set-TextCell ($ws.Range('D3')) '2.219.87'
Set-TextCell ($ws.Range('E4')) '  -0.09%  '
Set-TextCell ($ws.Range('D5')) '320.70'
Set-TextCell ($ws.Range('E5')) '  +0.10%  '
Set-TextCell ($ws.Range('D6')) '98.67'
Set-TextCell ($ws.Range('E6')) '  -9.34%  '
Set-TextCell ($ws.Range('D7')) '0.580'
Set-TextCell ($ws.Range('E7')) '  -9.12%  '
Set-TextCell ($ws.Range('E8')) '  -0.10%  '
Set-TextCell ($ws.Range('D9')) '0.563'
Set-TextCell ($ws.Range('E9')) '  -8.63%  '
Set-TextCell ($ws.Range('D10')) '36.77'
Set-TextCell ($ws.Range('E10')) '  -10.35%  '
Set-TextCell ($ws.Range('D11')) '54.02'
Set-TextCell ($ws.Range('E11')) '  -3.69%  '
Set-TextCell ($ws.Range('D12')) '0.0827'
Set-TextCell ($ws.Range('E12')) '  -10.25%  '
Set-TextCell ($ws.Range('D13')) '7.63'
Set-TextCell ($ws.Range('E13')) '  -10.44%  '
Set-TextCell ($ws.Range('D14')) '0.107'
Set-TextCell ($ws.Range('E14')) '  -2.21%  '
Set-TextCell ($ws.Range('D15')) '0.863'
Set-TextCell ($ws.Range('E15')) '  -12.15%  '
Set-TextCell ($ws.Range('D16')) '2.558.65'
Set-TextCell ($ws.Range('E16')) '  -6.54%  '
Set-TextCell ($ws.Range('E17')) '  -7.21%  '
Set-TextCell ($ws.Range('D18')) '2.219.83'
Set-TextCell ($ws.Range('E18')) '  -6.81%  '
Set-TextCell ($ws.Range('D19')) '42.946.60'
Set-TextCell ($ws.Range('E19')) '  -5.34%  '
Set-TextCell ($ws.Range('D20')) '13.99'
Set-TextCell ($ws.Range('E20')) '  -7.62%  '
Set-TextCell ($ws.Range('D21')) '0.0₃0964'
Set-TextCell ($ws.Range('E21')) '  -9.49%  '
Set-TextCell ($ws.Range('D22')) '6.54'
Set-TextCell ($ws.Range('E22')) '  -10.87%  '
Set-TextCell ($ws.Range('E23')) '  -12.71%  '
Set-TextCell ($ws.Range('D24')) '65.06'
Set-TextCell ($ws.Range('E24')) '  -11.28%  '
Set-TextCell ($ws.Range('D25')) '236.31'
Set-TextCell ($ws.Range('E25')) '  -10.63%  '
Set-TextCell ($ws.Range('D26')) '2.18'
Set-TextCell ($ws.Range('E26')) '  -6.94%  '
Set-TextCell ($ws.Range('E27')) '  -0.17%  '
Set-TextCell ($ws.Range('D28')) '4.04'
Set-TextCell ($ws.Range('E28')) '  +1.28%  '
Set-TextCell ($ws.Range('B29')) 'Toncoin'
Set-TextCell ($ws.Range('C29')) 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextCell ($ws.Range('D29')) '2.24'
Set-TextCell ($ws.Range('E29')) '  -2.74%  '
Set-TextCell ($ws.Range('B30')) 'Cosmos'
Set-TextCell ($ws.Range('C30')) 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextCell ($ws.Range('D30')) '9.97'
Set-TextCell ($ws.Range('E30')) '  -11.41%  '
Set-TextCell ($ws.Range('D31')) '6.34'
Set-TextCell ($ws.Range('D32')) '35.99'
Set-TextCell ($ws.Range('E32')) '  -3.45%  '
Set-TextCell ($ws.Range('D33')) '20.30'
Set-TextCell ($ws.Range('E33')) '  -9.56%  '
Set-TextCell ($ws.Range('D34')) '0.0864'
Set-TextCell ($ws.Range('E34')) '  -9.28%  '
Set-TextCell ($ws.Range('D35')) '153.45'
Set-TextCell ($ws.Range('E35')) '  -9.00%  '
Set-TextCell ($ws.Range('D36')) '2.67'
Set-TextCell ($ws.Range('E36')) '  -6.84%  '
Set-TextCell ($ws.Range('D37')) '3.24'
Set-TextCell ($ws.Range('E37')) '  -1.03%  '
Set-TextCell ($ws.Range('E38')) '  -7.75%  '
Set-TextCell ($ws.Range('D39')) '1.91'
Set-TextCell ($ws.Range('E39')) '  -2.35%  '
Set-TextCell ($ws.Range('D40')) '4.40'
Set-TextCell ($ws.Range('E40')) '  -6.96%  '
Set-TextCell ($ws.Range('E41')) '  -11.37%  '
Set-TextCell ($ws.Range('D42')) '3.71'
Set-TextCell ($ws.Range('E42')) '  -7.91%  '
Set-TextCell ($ws.Range('E43')) '  -9.54%  '
Set-TextCell ($ws.Range('D44')) '13.86'
Set-TextCell ($ws.Range('E44')) '  +6.72%  '
Set-TextCell ($ws.Range('D46')) '1.722.24'
Set-TextCell ($ws.Range('E46')) '  -8.35%  '
Set-TextCell ($ws.Range('B47')) 'BitcoinSV'
Set-TextCell ($ws.Range('C47')) 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
Set-TextCell ($ws.Range('D47')) '84.73'
Set-TextCell ($ws.Range('E47')) '  -13.90%  '
Set-TextCell ($ws.Range('B48')) 'Algorand'
Set-TextCell ($ws.Range('C48')) 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextCell ($ws.Range('D48')) '0.202'
Set-TextCell ($ws.Range('E48')) '  -11.73%  '
Set-TextCell ($ws.Range('E49')) '  -12.88%  '
Set-TextCell ($ws.Range('D50')) '8.80'
Set-TextCell ($ws.Range('E50')) '  -5.98%  '
Set-TextCell ($ws.Range('D51')) '74.73'
Set-TextCell ($ws.Range('E51')) '  -11.70%  '
